$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above row 99 (shifts the existing 99-188 block
# down to 101-190, which is exactly the "everything moves down by 2" pattern
# seen in the diff; the former last two rows end up as new rows 189-190).
$ws.Rows("99:100").Insert()

# Fill the two new rows (99 and 100) with the new weekly price records.
# Row 99: Primera
$ws.Cells.Item(99, 1).Value = 1
$ws.Cells.Item(99, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(99, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(99, 4).Value = 45240
$ws.Cells.Item(99, 5).Value = 15
$ws.Cells.Item(99, 6).Value = 100112042
$ws.Cells.Item(99, 7).Value = "Locoto"
$ws.Cells.Item(99, 8).Value = "Sin especificar"
$ws.Cells.Item(99, 9).Value = "Primera"
$ws.Cells.Item(99, 10).Value = 45
$ws.Cells.Item(99, 11).Value = 42000
$ws.Cells.Item(99, 12).Value = 45000
$ws.Cells.Item(99, 13).Value = 43000
$ws.Cells.Item(99, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(99, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(99, 16).Value = 2150
$ws.Cells.Item(99, 17).Value = 20
$ws.Cells.Item(99, 18).Value = "Hortaliza"

# Row 100: Segunda
$ws.Cells.Item(100, 1).Value = 1
$ws.Cells.Item(100, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(100, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(100, 4).Value = 45240
$ws.Cells.Item(100, 5).Value = 15
$ws.Cells.Item(100, 6).Value = 100112042
$ws.Cells.Item(100, 7).Value = "Locoto"
$ws.Cells.Item(100, 8).Value = "Sin especificar"
$ws.Cells.Item(100, 9).Value = "Segunda"
$ws.Cells.Item(100, 10).Value = 35
$ws.Cells.Item(100, 11).Value = 35000
$ws.Cells.Item(100, 12).Value = 38000
$ws.Cells.Item(100, 13).Value = 36286
$ws.Cells.Item(100, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(100, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(100, 16).Value = 1814
$ws.Cells.Item(100, 17).Value = 20
$ws.Cells.Item(100, 18).Value = "Hortaliza"
